# Rename the sheet from "Sheet1" to "Financial Model Building"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Financial Model Building"

# Re-apply the two formula blocks as multi-cell range formulas so Excel
# groups them into shared-formula blocks (matching the author's resave):
#  - D7:H8   -> master formula D$30*D35 (relative refs shift per cell)
#  - D15:H18 -> master formula D39      (relative refs shift per cell)
$ws.Range("D7:H8").Formula = "=D`$30*D35"
$ws.Range("D15:H18").Formula = "=D39"
